$wb = $excel.ActiveWorkbook

# Add the new "Branches" worksheet as the last (3rd) tab.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "Branches"

# Populate the branch/feature tree grid (write order matches the original authoring order).
$newSheet.Range('A1').Value = '1. Admin'
$newSheet.Range('B1').Value = '1.1 Parents'
$newSheet.Range('C1').Value = '1.2 Teachers'
$newSheet.Range('D1').Value = '1.3 Students'
$newSheet.Range('E1').Value = '1.4 Working Staff'
$newSheet.Range('B2').Value = '1.1.1 Childrens'
$newSheet.Range('D2').Value = '1.3.1 Time Table'
$newSheet.Range('D3').Value = '1.3.2 Events'
$newSheet.Range('C2').Value = '1.2.1 Classes assigned'
$newSheet.Range('C3').Value = '1.2.2 Subjects assigned'
$newSheet.Range('C4').Value = '1.2.3 Time table'
$newSheet.Range('C5').Value = '1.2.4 Events'
$newSheet.Range('C6').Value = '1.2.5 Notification(if any send by principal)'
$newSheet.Range('B3').Value = '1.1.2 Teacher assigned to students class'
$newSheet.Range('B4').Value = '1.1.3 Time table'
$newSheet.Range('B5').Value = '1.1.4 Subjects assigned'
$newSheet.Range('B6').Value = '1.1.5 Events'
$newSheet.Range('B7').Value = '1.1.6 Notification'
$newSheet.Range('B8').Value = '1.1.7 fees'
$newSheet.Range('D4').Value = '1.3.3 Notification'
$newSheet.Range('D5').Value = '1.3.4 Subjects'
$newSheet.Range('D6').Value = '1.3.5 Fees'
$newSheet.Range('F1').Value = '1.5 Departments'
$newSheet.Range('G1').Value = '1.6 Hostel'
$newSheet.Range('H1').Value = '1.7 Transportation'
$newSheet.Range('I1').Value = '1.8 Accounts'
$newSheet.Range('J1').Value = '1.9 Subjects'
$newSheet.Range('K1').Value = '1.10 Class'
$newSheet.Range('L1').Value = '1.11 Holidays'
$newSheet.Range('M1').Value = '1.12 Fees '
$newSheet.Range('N1').Value = '1.13 Library'
$newSheet.Range('O1').Value = '1.14 Exam '
$newSheet.Range('P1').Value = '1.15 Time Table'
$newSheet.Range('Q1').Value = '1.16 Events'
$newSheet.Range('R1').Value = '1.17 Printing'
$newSheet.Range('S1').Value = '1.18 Blog'
$newSheet.Range('C7').Value = '1.2.6 Attendance'
$newSheet.Range('D7').Value = '1.3.6 Attendance'
$newSheet.Range('B9').Value = '1.1.8 Attendance'
$newSheet.Range('T1').Value = '1.19 Notice'
$newSheet.Range('U1').Value = '1.20 Message'
$newSheet.Range('B10').Value = '1.1.9 Profile(self, childrens)'
$newSheet.Range('D8').Value = '1.3.7 Self Profile'
$newSheet.Range('D9').Value = '1.3.8 Exams, Exam result'
$newSheet.Range('B11').Value = '1.1.10 Exam, Exam result'
$newSheet.Range('C8').Value = '1.2.7 Salary'
$newSheet.Range('C9').Value = '1.2.8 Exams'
$newSheet.Range('E2').Value = '1.4.1 Salary'
$newSheet.Range('E3').Value = '1.4.2 Notification'
$newSheet.Range('G2').Value = '1.6.1 Room alloting'
$newSheet.Range('G3').Value = '1.6.2 Total student'
$newSheet.Range('G4').Value = '1.6.3 Fees Paid or unpaid'
$newSheet.Range('G5').Value = '1.6.4 Query module'
$newSheet.Range('D10').Value = '1.3.9 Query raise against hostel'
$newSheet.Range('I2').Value = '1.8.1 Slary of staff'
$newSheet.Range('I3').Value = '1.8.2 Fees of student'
$newSheet.Range('I4').Value = '1.8.3 Expenses'

# Match the author's final on-screen state: the whole table ends up selected
# on the new "Branches" tab (which becomes the active sheet/tab).
$newSheet.Activate()
$newSheet.Range('L1').Select() | Out-Null
$newSheet.Range('A1:U11').Activate() | Out-Null
